$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old "Localizer triggers" and "Explicit triggers" sections
# (rows 27-49), which will be rebuilt below with the new layout.
$ws.Range("A27:G49").Clear()

# --- Explicit triggers visual ---
$ws.Range("A28").Value = "Explicit triggers visual"
$ws.Range("A28").Font.Bold = $true
$ws.Range("A29").Value = "explicit_45_EXP"
$ws.Range("B29").Value = 121
$ws.Range("C29").Value = 125
$ws.Range("D29").Value = 129
$ws.Range("E29").Value = 133
$ws.Range("F29").Value = 137
$ws.Range("G29").Value = 141
$ws.Range("A30").Value = "explicit_45_UEX"
$ws.Range("B30").Value = 122
$ws.Range("C30").Value = 126
$ws.Range("D30").Value = 130
$ws.Range("E30").Value = 134
$ws.Range("F30").Value = 138
$ws.Range("G30").Value = 142
$ws.Range("A31").Value = "explicit_135_EXP"
$ws.Range("B31").Value = 123
$ws.Range("C31").Value = 127
$ws.Range("D31").Value = 131
$ws.Range("E31").Value = 135
$ws.Range("F31").Value = 139
$ws.Range("G31").Value = 143
$ws.Range("A32").Value = "explicit_135_UEX"
$ws.Range("B32").Value = 124
$ws.Range("C32").Value = 128
$ws.Range("D32").Value = 132
$ws.Range("E32").Value = 136
$ws.Range("F32").Value = 140
$ws.Range("G32").Value = 144
$ws.Range("A34").Value = "Explicit triggers auditory"
$ws.Range("A34").Font.Bold = $true
$ws.Range("A35").Value = "explicit_100_EXP"
$ws.Range("B35").Value = 145
$ws.Range("C35").Value = 149
$ws.Range("D35").Value = 153
$ws.Range("E35").Value = 157
$ws.Range("F35").Value = 161
$ws.Range("G35").Value = 165
$ws.Range("A36").Value = "explicit_100_UEX"
$ws.Range("B36").Value = 146
$ws.Range("C36").Value = 150
$ws.Range("D36").Value = 154
$ws.Range("E36").Value = 158
$ws.Range("F36").Value = 162
$ws.Range("G36").Value = 166
$ws.Range("A37").Value = "explicit_160_EXP"
$ws.Range("B37").Value = 147
$ws.Range("C37").Value = 151
$ws.Range("D37").Value = 155
$ws.Range("E37").Value = 159
$ws.Range("F37").Value = 163
$ws.Range("G37").Value = 167
$ws.Range("A38").Value = "explicit_160_UEX"
$ws.Range("B38").Value = 148
$ws.Range("C38").Value = 152
$ws.Range("D38").Value = 156
$ws.Range("E38").Value = 160
$ws.Range("F38").Value = 164
$ws.Range("G38").Value = 168
$ws.Range("A40").Value = "Localizer triggers"
$ws.Range("A40").Font.Bold = $true
$ws.Range("A41").Value = "loc_start"
$ws.Range("B41").Value = 169
$ws.Range("A42").Value = "loc_isi"
$ws.Range("B42").Value = 170
$ws.Range("A43").Value = "loc_45_100"
$ws.Range("B43").Value = 171
$ws.Range("A44").Value = "loc_45_160"
$ws.Range("B44").Value = 172
$ws.Range("A45").Value = "loc_135_100"
$ws.Range("B45").Value = 173
$ws.Range("A46").Value = "loc_135_160"
$ws.Range("B46").Value = 174
$ws.Range("A47").Value = "loc_response"
$ws.Range("B47").Value = 175

# Move the active selection to match the saved view state in the target
# workbook.
$ws.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H43").Select()
